$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 28
$ws1.Range("F4").Value = 226
$ws1.Range("F5").Value = 3811
$ws1.Range("F7").Value = 22
$ws1.Range("F8").Value = 434

# Sheet "全部类型" (All types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 28
$ws4.Range("F4").Value = 226
$ws4.Range("F5").Value = 3811
$ws4.Range("F9").Value = 22
$ws4.Range("F10").Value = 434
